# "show current date and time"
# Adds a small "current date/time" demo block to the FUNCTIONS section
# (columns R/T) of the BUDGET sheet: System Date (NOW), Date (TODAY),
# Day, Month and Year (DAY/MONTH/YEAR of TODAY()).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BUDGET")

# Row 26 - System Date (current date & time)
$ws.Range("R26").Value = "System Date"
$ws.Range("T26").NumberFormat = "m/d/yy h:mm"
$ws.Range("T26").Formula = "=NOW()"

# Row 27 - Date (today, date only)
$ws.Range("R27").Value = "Date"
$ws.Range("T27").NumberFormat = "mm-dd-yy"
$ws.Range("T27").Formula = "=TODAY()"

# Row 28 - Day component of today
$ws.Range("R28").Value = "Day"
$ws.Range("T28").Formula = "=DAY(TODAY())"

# Row 29 - Month component of today
$ws.Range("R29").Value = "Month"
$ws.Range("T29").Formula = "=MONTH(TODAY())"

# Row 30 - Year component of today
$ws.Range("R30").Value = "Year"
$ws.Range("T30").Formula = "=YEAR(TODAY())"

# Scroll the window down toward the new rows and leave the selection where
# the author ended up (R34), matching the saved view state.
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 7
$ws.Range("R34").Select()
